$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 3
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 1.428571428571429
$ws.Cells.Item(2, 7).Value = 5
$ws.Cells.Item(2, 8).Value = 1.666666666666667
$ws.Cells.Item(2, 9).Value = 15
$ws.Cells.Item(2, 10).Value = 2.142857142857143
$ws.Cells.Item(2, 11).Value = 4
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2
$ws.Cells.Item(2, 14).Value = 0.6666666666666666
$ws.Cells.Item(2, 15).Value = 6
$ws.Cells.Item(2, 16).Value = 0.8571428571428571

$ws.Cells.Item(3, 2).Value = 7
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(3, 6).Value = 0.8571428571428571
$ws.Cells.Item(3, 7).Value = 6
$ws.Cells.Item(3, 8).Value = 2
$ws.Cells.Item(3, 9).Value = 12
$ws.Cells.Item(3, 10).Value = 1.714285714285714
$ws.Cells.Item(3, 11).Value = 5
$ws.Cells.Item(3, 12).Value = 1.25
$ws.Cells.Item(3, 13).Value = 4
$ws.Cells.Item(3, 14).Value = 1.333333333333333
$ws.Cells.Item(3, 15).Value = 9
$ws.Cells.Item(3, 16).Value = 1.285714285714286

$ws.Cells.Item(4, 2).Value = 7
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 4
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 0.5714285714285714
$ws.Cells.Item(4, 7).Value = 4
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 10).Value = 1.142857142857143
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7
$ws.Cells.Item(4, 14).Value = 1.75
$ws.Cells.Item(4, 15).Value = 10
$ws.Cells.Item(4, 16).Value = 1.428571428571429

$ws.Cells.Item(5, 2).Value = 7
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = 3
$ws.Cells.Item(5, 5).Value = 11
$ws.Cells.Item(5, 6).Value = 1.571428571428571
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = 0.6666666666666666
$ws.Cells.Item(5, 9).Value = 13
$ws.Cells.Item(5, 10).Value = 1.857142857142857
$ws.Cells.Item(5, 11).Value = 6
$ws.Cells.Item(5, 12).Value = 1.5
$ws.Cells.Item(5, 13).Value = 7
$ws.Cells.Item(5, 14).Value = 2.333333333333333
$ws.Cells.Item(5, 15).Value = 13
$ws.Cells.Item(5, 16).Value = 1.857142857142857

$ws.Cells.Item(6, 2).Value = 7
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6
$ws.Cells.Item(6, 8).Value = 2
$ws.Cells.Item(6, 9).Value = 13
$ws.Cells.Item(6, 10).Value = 1.857142857142857
$ws.Cells.Item(6, 11).Value = 5
$ws.Cells.Item(6, 12).Value = 1.25
$ws.Cells.Item(6, 13).Value = 5
$ws.Cells.Item(6, 14).Value = 1.666666666666667
$ws.Cells.Item(6, 15).Value = 10
$ws.Cells.Item(6, 16).Value = 1.428571428571429

$ws.Cells.Item(7, 2).Value = 7
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 0.8571428571428571
$ws.Cells.Item(7, 7).Value = 10
$ws.Cells.Item(7, 8).Value = 3.333333333333333
$ws.Cells.Item(7, 9).Value = 16
$ws.Cells.Item(7, 10).Value = 2.285714285714286
$ws.Cells.Item(7, 11).Value = 6
$ws.Cells.Item(7, 12).Value = 1.5
$ws.Cells.Item(7, 13).Value = 2
$ws.Cells.Item(7, 14).Value = 0.6666666666666666
$ws.Cells.Item(7, 15).Value = 8
$ws.Cells.Item(7, 16).Value = 1.142857142857143

$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 3).Value = 4
$ws.Cells.Item(8, 4).Value = 3
$ws.Cells.Item(8, 5).Value = 11
$ws.Cells.Item(8, 6).Value = 1.571428571428571
$ws.Cells.Item(8, 7).Value = 6
$ws.Cells.Item(8, 8).Value = 2
$ws.Cells.Item(8, 9).Value = 17
$ws.Cells.Item(8, 10).Value = 2.428571428571428
$ws.Cells.Item(8, 11).Value = 6
$ws.Cells.Item(8, 12).Value = 1.5
$ws.Cells.Item(8, 13).Value = 2
$ws.Cells.Item(8, 14).Value = 0.6666666666666666
$ws.Cells.Item(8, 15).Value = 8
$ws.Cells.Item(8, 16).Value = 1.142857142857143

$ws.Cells.Item(9, 2).Value = 7
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.2857142857142857
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 5
$ws.Cells.Item(9, 10).Value = 0.7142857142857143
$ws.Cells.Item(9, 11).Value = 5
$ws.Cells.Item(9, 12).Value = 1.25
$ws.Cells.Item(9, 13).Value = 5
$ws.Cells.Item(9, 14).Value = 1.666666666666667
$ws.Cells.Item(9, 15).Value = 10
$ws.Cells.Item(9, 16).Value = 1.428571428571429

$ws.Cells.Item(10, 2).Value = 7
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 0.5714285714285714
$ws.Cells.Item(10, 7).Value = 3
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 7
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 7
$ws.Cells.Item(10, 12).Value = 1.75
$ws.Cells.Item(10, 13).Value = 8
$ws.Cells.Item(10, 14).Value = 2.666666666666667
$ws.Cells.Item(10, 15).Value = 15
$ws.Cells.Item(10, 16).Value = 2.142857142857143

$ws.Cells.Item(11, 2).Value = 7
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 4
$ws.Cells.Item(11, 5).Value = 6
$ws.Cells.Item(11, 6).Value = 0.8571428571428571
$ws.Cells.Item(11, 7).Value = 4
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = 10
$ws.Cells.Item(11, 10).Value = 1.428571428571429
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5
$ws.Cells.Item(11, 14).Value = 1.25
$ws.Cells.Item(11, 15).Value = 8
$ws.Cells.Item(11, 16).Value = 1.142857142857143

$ws.Cells.Item(12, 2).Value = 7
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 4
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 0.4285714285714285
$ws.Cells.Item(12, 7).Value = 3
$ws.Cells.Item(12, 8).Value = 0.75
$ws.Cells.Item(12, 9).Value = 6
$ws.Cells.Item(12, 10).Value = 0.8571428571428571
$ws.Cells.Item(12, 11).Value = 5
$ws.Cells.Item(12, 12).Value = 1.666666666666667
$ws.Cells.Item(12, 13).Value = 9
$ws.Cells.Item(12, 14).Value = 2.25
$ws.Cells.Item(12, 15).Value = 14
$ws.Cells.Item(12, 16).Value = 2

$ws.Cells.Item(13, 2).Value = 7
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(13, 4).Value = 3
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 0.5714285714285714
$ws.Cells.Item(13, 7).Value = 5
$ws.Cells.Item(13, 8).Value = 1.666666666666667
$ws.Cells.Item(13, 9).Value = 9
$ws.Cells.Item(13, 10).Value = 1.285714285714286
$ws.Cells.Item(13, 11).Value = 4
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 8
$ws.Cells.Item(13, 14).Value = 2.666666666666667
$ws.Cells.Item(13, 15).Value = 12
$ws.Cells.Item(13, 16).Value = 1.714285714285714

$ws.Cells.Item(14, 2).Value = 7
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 4
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = 0.7142857142857143
$ws.Cells.Item(14, 7).Value = 8
$ws.Cells.Item(14, 8).Value = 2
$ws.Cells.Item(14, 9).Value = 13
$ws.Cells.Item(14, 10).Value = 1.857142857142857
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 1
$ws.Cells.Item(14, 14).Value = 0.25
$ws.Cells.Item(14, 15).Value = 2
$ws.Cells.Item(14, 16).Value = 0.2857142857142857

$ws.Cells.Item(15, 2).Value = 7
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 4
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.1428571428571428
$ws.Cells.Item(15, 7).Value = 4
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 5
$ws.Cells.Item(15, 10).Value = 0.7142857142857143
$ws.Cells.Item(15, 11).Value = 6
$ws.Cells.Item(15, 12).Value = 2
$ws.Cells.Item(15, 13).Value = 2
$ws.Cells.Item(15, 14).Value = 0.5
$ws.Cells.Item(15, 15).Value = 8
$ws.Cells.Item(15, 16).Value = 1.142857142857143

$ws.Cells.Item(16, 2).Value = 7
$ws.Cells.Item(16, 3).Value = 3
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 0.5714285714285714
$ws.Cells.Item(16, 7).Value = 4
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 8
$ws.Cells.Item(16, 10).Value = 1.142857142857143
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 5
$ws.Cells.Item(16, 14).Value = 1.25
$ws.Cells.Item(16, 15).Value = 7
$ws.Cells.Item(16, 16).Value = 1

$ws.Cells.Item(17, 2).Value = 7
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 4).Value = 4
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.2857142857142857
$ws.Cells.Item(17, 7).Value = 5
$ws.Cells.Item(17, 8).Value = 1.25
$ws.Cells.Item(17, 9).Value = 7
$ws.Cells.Item(17, 10).Value = 1
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3
$ws.Cells.Item(17, 14).Value = 0.75
$ws.Cells.Item(17, 15).Value = 6
$ws.Cells.Item(17, 16).Value = 0.8571428571428571

$ws.Cells.Item(18, 2).Value = 7
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.1428571428571428
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 0.75
$ws.Cells.Item(18, 9).Value = 4
$ws.Cells.Item(18, 10).Value = 0.5714285714285714
$ws.Cells.Item(18, 11).Value = 5
$ws.Cells.Item(18, 12).Value = 1.666666666666667
$ws.Cells.Item(18, 13).Value = 10
$ws.Cells.Item(18, 14).Value = 2.5
$ws.Cells.Item(18, 15).Value = 15
$ws.Cells.Item(18, 16).Value = 2.142857142857143

$ws.Cells.Item(19, 2).Value = 7
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = 4
$ws.Cells.Item(19, 5).Value = 7
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 7
$ws.Cells.Item(19, 8).Value = 1.75
$ws.Cells.Item(19, 9).Value = 14
$ws.Cells.Item(19, 10).Value = 2
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 0.6666666666666666
$ws.Cells.Item(19, 13).Value = 6
$ws.Cells.Item(19, 14).Value = 1.5
$ws.Cells.Item(19, 15).Value = 8
$ws.Cells.Item(19, 16).Value = 1.142857142857143

$ws.Cells.Item(20, 2).Value = 7
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = 6
$ws.Cells.Item(20, 6).Value = 0.8571428571428571
$ws.Cells.Item(20, 7).Value = 4
$ws.Cells.Item(20, 8).Value = 1.333333333333333
$ws.Cells.Item(20, 9).Value = 10
$ws.Cells.Item(20, 10).Value = 1.428571428571429
$ws.Cells.Item(20, 11).Value = 9
$ws.Cells.Item(20, 12).Value = 2.25
$ws.Cells.Item(20, 13).Value = 2
$ws.Cells.Item(20, 14).Value = 0.6666666666666666
$ws.Cells.Item(20, 15).Value = 11
$ws.Cells.Item(20, 16).Value = 1.571428571428571

$ws.Cells.Item(21, 2).Value = 7
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(21, 4).Value = 4
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 0.5714285714285714
$ws.Cells.Item(21, 7).Value = 5
$ws.Cells.Item(21, 8).Value = 1.25
$ws.Cells.Item(21, 9).Value = 9
$ws.Cells.Item(21, 10).Value = 1.285714285714286
$ws.Cells.Item(21, 11).Value = 10
$ws.Cells.Item(21, 12).Value = 3.333333333333333
$ws.Cells.Item(21, 13).Value = 11
$ws.Cells.Item(21, 14).Value = 2.75
$ws.Cells.Item(21, 15).Value = 21
$ws.Cells.Item(21, 16).Value = 3

